$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# New log entries for 2016-09-27 (serial 42640), appended after the
# existing last row (565). Rows 566-569 are intentionally left blank,
# matching the source data gap.

# Row 570
$ws.Range("A570").Value2 = "Demo"
$ws.Range("B570").Value2 = 42640
$ws.Range("C570").Value2 = "1600"
$ws.Range("D570").Value2 = "DB"
$ws.Range("E570").Value2 = "0011"

# Row 571
$ws.Range("A571").Value2 = "Demo"
$ws.Range("B571").Value2 = 42640
$ws.Range("C571").Value2 = "1600"
$ws.Range("D571").Value2 = "DB"
$ws.Range("E571").Value2 = "0004"

# Row 572
$ws.Range("A572").Value2 = "Demo"
$ws.Range("B572").Value2 = 42640
$ws.Range("C572").Value2 = "1630"
$ws.Range("D572").Value2 = "OSG"
$ws.Range("E572").Value2 = "2001"

# Row 573
$ws.Range("A573").Value2 = "Operator"
$ws.Range("B573").Value2 = 42640
$ws.Range("C573").Value2 = "1530"
$ws.Range("D573").Value2 = "KT"
$ws.Range("E573").Value2 = "524"
$ws.Range("F573").Value2 = "VIP - please replace Alex at 3:30"

# Row 574
$ws.Range("A574").Value2 = "Operator"
$ws.Range("B574").Value2 = 42640
$ws.Range("C574").Value2 = "1830"
$ws.Range("D574").Value2 = "OSG"
$ws.Range("E574").Value2 = "1001"
$ws.Range("F574").Value2 = "Please remain on site and oversee rooms 1001, 1005, 1006, 2003, 2004 recordings"
$ws.Rows.Item(574).RowHeight = 30

# Row 575
$ws.Range("A575").Value2 = "Setup Skype Kit"
$ws.Range("B575").Value2 = 42640
$ws.Range("C575").Value2 = "1830"
$ws.Range("D575").Value2 = "OSG"
$ws.Range("E575").Value2 = "1005"
$ws.Range("F575").Value2 = "Video recording via WinMovie  maker -  web cam and tripod in OSG 1014L"
$ws.Rows.Item(575).RowHeight = 30

# Row 576
$ws.Range("A576").Value2 = "Setup Skype Kit"
$ws.Range("B576").Value2 = 42640
$ws.Range("C576").Value2 = "1830"
$ws.Range("D576").Value2 = "OSG"
$ws.Range("E576").Value2 = "1006"
$ws.Range("F576").Value2 = "Video recording via WinMovie  maker -  web cam and tripod in OSG 1014L"
$ws.Rows.Item(576).RowHeight = 30

# Row 577
$ws.Range("A577").Value2 = "Setup Skype Kit"
$ws.Range("B577").Value2 = 42640
$ws.Range("C577").Value2 = "1830"
$ws.Range("D577").Value2 = "OSG"
$ws.Range("E577").Value2 = "2003"
$ws.Range("F577").Value2 = "Video recording via WinMovie  maker -  web cam and tripod in OSG 1014L"
$ws.Rows.Item(577).RowHeight = 30

# Row 578
$ws.Range("A578").Value2 = "Setup Skype Kit"
$ws.Range("B578").Value2 = 42640
$ws.Range("C578").Value2 = "1830"
$ws.Range("D578").Value2 = "OSG"
$ws.Range("E578").Value2 = "2004"
$ws.Range("F578").Value2 = "Video recording via WinMovie  maker -  web cam and tripod in OSG 1014L"
$ws.Rows.Item(578).RowHeight = 30

# Row 579
$ws.Range("A579").Value2 = "Setup Skype Kit"
$ws.Range("B579").Value2 = 42640
$ws.Range("C579").Value2 = "1830"
$ws.Range("D579").Value2 = "OSG"
$ws.Range("E579").Value2 = "1001"
$ws.Range("F579").Value2 = "Video recording via WinMovie  maker -  web cam and tripod in OSG 1014L"
$ws.Rows.Item(579).RowHeight = 30

# Row 580
$ws.Range("A580").Value2 = "Pickup Skype Kit"
$ws.Range("B580").Value2 = 42640
$ws.Range("C580").Value2 = "2150"
$ws.Range("D580").Value2 = "OSG"
$ws.Range("E580").Value2 = "1001"
$ws.Range("F580").Value2 = "Return web cam and tripod to OSG 1014L"

# Row 581
$ws.Range("A581").Value2 = "Pickup Skype Kit"
$ws.Range("B581").Value2 = 42640
$ws.Range("C581").Value2 = "2150"
$ws.Range("D581").Value2 = "OSG"
$ws.Range("E581").Value2 = "1005"
$ws.Range("F581").Value2 = "Return web cam and tripod to OSG 1014L"

# Row 582
$ws.Range("A582").Value2 = "Pickup Skype Kit"
$ws.Range("B582").Value2 = 42640
$ws.Range("C582").Value2 = "2150"
$ws.Range("D582").Value2 = "OSG"
$ws.Range("E582").Value2 = "1006"
$ws.Range("F582").Value2 = "Return web cam and tripod to OSG 1014L"

# Row 583
$ws.Range("A583").Value2 = "Pickup Skype Kit"
$ws.Range("B583").Value2 = 42640
$ws.Range("C583").Value2 = "2150"
$ws.Range("D583").Value2 = "OSG"
$ws.Range("E583").Value2 = "2003"
$ws.Range("F583").Value2 = "Return web cam and tripod to OSG 1014L"

# Row 584
$ws.Range("A584").Value2 = "Pickup Skype Kit"
$ws.Range("B584").Value2 = 42640
$ws.Range("C584").Value2 = "2150"
$ws.Range("D584").Value2 = "OSG"
$ws.Range("E584").Value2 = "2004"
$ws.Range("F584").Value2 = "Return web cam and tripod to OSG 1014L"

# Update the view: scroll so row 566 is at the top and the next empty
# cell (F587) is selected, ready for data entry.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 566
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F587").Select()
